$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelRuns")

# 1. Mark the 2005 "v2" run (row 4) as the current run by adding "current" in column H.
$ws.Range("H4").Value = "current"

# 2. Insert a new row before the old row 76 (2050 row) to hold the new
#    2035 IPA 11 run, shifting the old row 76 (and everything after it) down to row 77.
$ws.Rows("76:76").Insert(1) # xlShiftDown

# Apply the same formatting as the neighboring 2035 row (row 73), which uses the
# plain style (s="12") in every column, including O/P (unlike row 75, whose O/P
# cells are blank and use a different fill style).
$ws.Range("A73:S73").Copy()
$ws.Range("A76:S76").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Populate the new row with the 2035 IPA 11 run data.
$ws.Range("A76").Value = 2035
$ws.Range("B76").Value = "2035_TM160_IPA_11"
$ws.Range("C76").Value = "RTP2025_IP"
$ws.Range("D76").Value = "IPA"
$ws.Range("E76").Value = "IPA with EN7 fixed and bike mode share adjusted"
$ws.Range("F76").Value = "FBP scaled to RGF"
$ws.Range("G76").Value = "run182"
$ws.Range("H76").Value = "current"
$ws.Range("I76").Value = "M:\Application\Model One\RTP2021\Blueprint\INPUT_DEVELOPMENT\Networks\BlueprintNetworks_64\net_2035_Blueprint_tollscsv"
$ws.Range("J76").Value = "model3-c"
$ws.Range("K76").Value = "https://app.asana.com/0/1204085012544660/1206021318810361/f"
$ws.Range("L76").Value = 20.55
$ws.Range("M76").Value = "na"
$ws.Range("N76").Value = "na"
$ws.Range("O76").Value = 0.87
$ws.Range("P76").Value = 0.78
$ws.Range("Q76").Value = 100
$ws.Range("R76").Value = 0
$ws.Range("S76").Value = 75

$wb.Save()
